$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 7 is a re-add of the original row 2 data (Svetlana/Marhefka/sm269@njit.edu),
#    so copy row 2 -> row 7 first (values + formatting) before row 2 itself is edited.
$ws.Range("A2:E2").Copy($ws.Range("A7:E7"))

# 2) New row 6: merged "Total" label across A6:B6, centered.
$ws.Range("A6:B6").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A6").Value = "Total"
$ws.Range("A6:B6").Merge()

# 3) Row 2 edits: first name Svetlana -> Lana, favorite number 13 -> 66.
$ws.Range("A2").Value = "Lana"
$ws.Range("D2").Value = 66

# 4) Row 7's favorite number differs from the original row 2 value (13 -> 77).
$ws.Range("D7").Value = 77

# 5) New hyperlink on C7 (mirrors the one on C2); re-apply C2's formatting afterward
#    since adding a hyperlink resets the cell's style.
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:sm269@njit.edu")
$ws.Range("C7").Style = $ws.Range("C2").Style

# 6) Active selection is on A3.
$ws.Range("A3").Select()
